$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 243
$ws.Range("I2").Value = 604
$ws.Range("J2").Value = 2520
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 738
$ws.Range("M2").Value = 33
$ws.Range("N2").Value = 409
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 23
$ws.Range("S2").Value = 276
$ws.Range("T2").Value = 428
$ws.Range("U2").Value = 48
$ws.Range("V2").Value = 3963
$ws.Range("X2").Value = 3950
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 61
$ws.Range("AA2").Value = 26
